$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows before row 8 (pushes "Analytic" row 8 -> 10,
#    "Number of possible errors" row 10 -> 12, "1 bit Hb(p)" row 12 -> 14,
#    and the entropy-calc row 13 -> 15). Formulas referencing A13 auto-update
#    to A15 because this is a real row insert.
# ---------------------------------------------------------------------------
$ws.Rows("8:9").Insert()

# New iteration-count rows for the larger lattice sizes.
$ws.Range("A8").Value = 100000000
$ws.Range("A9").Value = 1000000000

# The insert carries formatting down from row 7 into the blank new rows 8:9
# for columns B:Z - strip that back off so those cells stay empty/unstyled.
$ws.Range("B8:Z9").ClearFormats()

# ---------------------------------------------------------------------------
# 2. Fill in the previously-empty 10^7 totals (row 7) for the 5th sub-column
#    of each L block (H(S,EC) grouping reaching its final steady value).
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = 22.928278335105901
$ws.Range("K7").Value = 22.820774259042
$ws.Range("P7").Value = 0.83716017215875405
$ws.Range("U7").Value = 22.8337258456556
$ws.Range("Z7").Value = 0.82420858554516396

# ---------------------------------------------------------------------------
# 3. Style cleanup across the header row and the five data rows (rows 2-7):
#    the thin inner-column cells (C:E, H:J, M:O, R:T, W:Y) drop their border
#    style entirely, and the section-divider column (F, K, P, U) on the
#    header row matches the plain divider style already used below it.
# ---------------------------------------------------------------------------
$ws.Range("C2:E7").ClearFormats()
$ws.Range("H2:J7").ClearFormats()
$ws.Range("M2:O7").ClearFormats()
$ws.Range("R2:T7").ClearFormats()
$ws.Range("W2:Y7").ClearFormats()

$ws.Range("F3").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("P2").PasteSpecial(-4122)
$ws.Range("U2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Move the convergence chart further down the sheet, below the data that
#    now extends two rows lower (same size, new anchor position).
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Top = 349.3285826771654
$co.Left = 48.34291338582677
$co.Width = 521.0673828125
$co.Height = 216

# ---------------------------------------------------------------------------
# 5. View state: zoom to 70% and move the active selection.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("N22").Select()
